$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 47: BATCH046 gets aborted ---
$ws.Range("F47").Value = "Aborted"
$ws.Range("N47").Value = "no need"
$ws.Range("O47").Value = "adibsv"
$ws.Range("P47").Value = "2025-11-05 11:57:01"

# --- Row 48: new batch BATCH047 ---
$ws.Range("A48").Value = "BATCH047"
$ws.Range("B48").Value = "RCP001"
$ws.Range("C48").Value = "Official - Recipe 01"
$ws.Range("D48").Value = 2
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = "InProgress"
$ws.Range("G48").Value = "2025-11-06 10:55:00"
$ws.Range("H48").Value = "2025-11-06 11:55:00"
$ws.Range("I48").Value = "adibsv"
$ws.Range("J48").Value = 45966.4969015509
$ws.Range("J48").NumberFormat = "m/d/yy h:mm"
$ws.Range("K48").Value = "adibsv"
$ws.Range("L48").Value = "2025-11-05 11:55:40"
$ws.Range("M48").Value = ""
$ws.Range("N48").Value = ""
$ws.Range("O48").Value = ""
$ws.Range("P48").Value = ""
$ws.Range("Q48").Value = ""

# --- Row 49: new batch BATCH048 ---
$ws.Range("A49").Value = "BATCH048"
$ws.Range("B49").Value = "RCP001"
$ws.Range("C49").Value = "Official - Recipe 01"
$ws.Range("D49").Value = 1
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = "InProgress"
$ws.Range("G49").Value = "2025-11-06 12:10:00"
$ws.Range("H49").Value = "2025-11-08 12:10:00"
$ws.Range("I49").Value = "adibsv"
$ws.Range("J49").Value = 45966.507346088
$ws.Range("J49").NumberFormat = "m/d/yy h:mm"
$ws.Range("K49").Value = "adibsv"
$ws.Range("L49").Value = "2025-11-05 12:10:36"
$ws.Range("M49").Value = ""
$ws.Range("N49").Value = ""
$ws.Range("O49").Value = ""
$ws.Range("P49").Value = ""
$ws.Range("Q49").Value = ""
